$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Habilidades sheet (sheet7): insert a new "grupo" column (C) that groups
# skills, shifting the former C/D ("nome_pt"/"nome_en") columns to D/E.
# ---------------------------------------------------------------------------
$wsHab = $wb.Worksheets.Item("Habilidades")

$wsHab.Columns("C").Insert() | Out-Null

# Give the new column the same custom width as column B.
$wsHab.Columns("C").ColumnWidth = $wsHab.Columns("B").ColumnWidth

# Header + first data row for the new "grupo" column.
$wsHab.Range("C1").Value = "grupo"
$wsHab.Range("C2").Value = "a"

# ---------------------------------------------------------------------------
# Classes sheet (sheet8): fill in the missing "comportamental" alias cell
# for the behavioural-type row (row 5), matching the already-present value
# in column E, with an underlined style.
# ---------------------------------------------------------------------------
$wsClasses = $wb.Worksheets.Item("Classes")
$wsClasses.Range("B5").Value = "comportamental"
$wsClasses.Range("B5").Font.Underline = 2

# ---------------------------------------------------------------------------
# Restore per-sheet selections (cursor position the author left each sheet
# in before saving).
# ---------------------------------------------------------------------------
$wsCabecalho = $wb.Worksheets.Item("Cabeçalho")
$wsCabecalho.Activate() | Out-Null
$wsCabecalho.Range("B1:C1").Select() | Out-Null

$wsResumo = $wb.Worksheets.Item("Resumo")
$wsResumo.Activate() | Out-Null
$wsResumo.Range("D31").Select() | Out-Null

$wsOutros = $wb.Worksheets.Item("Outros")
$wsOutros.Activate() | Out-Null
$wsOutros.Range("A1:E1").Select() | Out-Null

$wsExperiencias = $wb.Worksheets.Item("Experiências")
$wsExperiencias.Activate() | Out-Null
$wsExperiencias.Rows("1:1").Select() | Out-Null

$wsFormacoes = $wb.Worksheets.Item("Formações")
$wsFormacoes.Activate() | Out-Null
$wsFormacoes.Range("C10").Select() | Out-Null

$wsClasses.Activate() | Out-Null
$wsClasses.Range("B5").Select() | Out-Null

# Habilidades ends up the active tab/sheet, matching the saved workbook view.
$wsHab.Activate() | Out-Null
$wsHab.Range("G8").Select() | Out-Null
